# Scheduled market-data refresh: update leve profit calculator inputs
# (currentAveragePrice/NQ/HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ) across
# all job sheets to the latest pulled prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3946.7568
$ws.Range("J112").Value = 3998.0557
$ws.Range("L112").Value = 11994.1671
$ws.Range("N112").Value = -14210.1671
$ws.Range("H129").Value = 1701.5454
$ws.Range("I129").Value = 664.625
$ws.Range("J129").Value = 4466.6665
$ws.Range("K129").Value = 1993.875
$ws.Range("L129").Value = 13399.9995
$ws.Range("M129").Value = 3006.125
$ws.Range("N129").Value = -23399.9995
$ws.Range("H137").Value = 1184.9412
$ws.Range("I137").Value = 548.9
$ws.Range("J137").Value = 2093.5715
$ws.Range("K137").Value = 1646.7
$ws.Range("L137").Value = 6280.7145
$ws.Range("M137").Value = 903.3000000000002
$ws.Range("N137").Value = -11380.7145
$ws.Range("H138").Value = 3172.0408
$ws.Range("I138").Value = 3527.2
$ws.Range("J138").Value = 3080.9744
$ws.Range("K138").Value = 10581.6
$ws.Range("L138").Value = 9242.923200000001
$ws.Range("M138").Value = -5441.599999999999
$ws.Range("N138").Value = -19522.9232
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5954135
$ws.Range("I61").Value = 9260792
$ws.Range("J61").Value = 2151.85
$ws.Range("K61").Value = 9260792
$ws.Range("L61").Value = 2151.85
$ws.Range("M61").Value = -9260580
$ws.Range("N61").Value = -2575.85
$ws.Range("H74").Value = 1314.7037
$ws.Range("I74").Value = 1550.75
$ws.Range("J74").Value = 1125.8667
$ws.Range("K74").Value = 1550.75
$ws.Range("L74").Value = 1125.8667
$ws.Range("M74").Value = -676.75
$ws.Range("N74").Value = -2873.8667
$ws.Range("H77").Value = 1314.7037
$ws.Range("I77").Value = 1550.75
$ws.Range("J77").Value = 1125.8667
$ws.Range("K77").Value = 7753.75
$ws.Range("L77").Value = 5629.333500000001
$ws.Range("M77").Value = -3385.75
$ws.Range("N77").Value = -14365.3335
$ws.Range("H136").Value = 5954135
$ws.Range("I136").Value = 9260792
$ws.Range("J136").Value = 2151.85
$ws.Range("K136").Value = 27782376
$ws.Range("L136").Value = 6455.549999999999
$ws.Range("M136").Value = -27779826
$ws.Range("N136").Value = -11555.55
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 75200.5
$ws.Range("J92").Value = 75200.5
$ws.Range("L92").Value = 75200.5
$ws.Range("N92").Value = -80192.5
$ws.Range("H134").Value = 2918.2856
$ws.Range("I134").Value = 2967.75
$ws.Range("J134").Value = 2810.3635
$ws.Range("K134").Value = 8903.25
$ws.Range("L134").Value = 8431.0905
$ws.Range("M134").Value = -6368.25
$ws.Range("N134").Value = -13501.0905
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1516.5
$ws.Range("I31").Value = 1337.3334
$ws.Range("J31").Value = 1785.25
$ws.Range("K31").Value = 1337.3334
$ws.Range("L31").Value = 1785.25
$ws.Range("M31").Value = -1042.3334
$ws.Range("N31").Value = -2375.25
$ws.Range("H34").Value = 1516.5
$ws.Range("I34").Value = 1337.3334
$ws.Range("J34").Value = 1785.25
$ws.Range("K34").Value = 1337.3334
$ws.Range("L34").Value = 1785.25
$ws.Range("M34").Value = -1135.3334
$ws.Range("N34").Value = -2189.25
$ws.Range("H58").Value = 1645.7567
$ws.Range("I58").Value = 1206.5264
$ws.Range("J58").Value = 2109.389
$ws.Range("K58").Value = 1206.5264
$ws.Range("L58").Value = 2109.389
$ws.Range("M58").Value = -1003.5264
$ws.Range("N58").Value = -2515.389
$ws.Range("H136").Value = 1645.7567
$ws.Range("I136").Value = 1206.5264
$ws.Range("J136").Value = 2109.389
$ws.Range("K136").Value = 3619.5792
$ws.Range("L136").Value = 6328.167
$ws.Range("M136").Value = -1069.5792
$ws.Range("N136").Value = -11428.167
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 921.2
$ws.Range("J68").Value = 986.75
$ws.Range("L68").Value = 2960.25
$ws.Range("N68").Value = -4582.25
$ws.Range("H71").Value = 921.2
$ws.Range("J71").Value = 986.75
$ws.Range("L71").Value = 8880.75
$ws.Range("N71").Value = -16992.75
$ws.Range("H107").Value = 1403.4828
$ws.Range("I107").Value = 268.08
$ws.Range("J107").Value = 2263.6365
$ws.Range("K107").Value = 804.24
$ws.Range("L107").Value = 6790.9095
$ws.Range("M107").Value = 1115.76
$ws.Range("N107").Value = -10630.9095
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1251.7059
$ws.Range("I122").Value = 1194.4546
$ws.Range("J122").Value = 1356.6666
$ws.Range("K122").Value = 3583.3638
$ws.Range("L122").Value = 4069.9998
$ws.Range("M122").Value = -1133.3638
$ws.Range("N122").Value = -8969.9998
$ws.Range("H132").Value = 2566.923
$ws.Range("I132").Value = 2811.2727
$ws.Range("J132").Value = 2387.7334
$ws.Range("K132").Value = 8433.8181
$ws.Range("L132").Value = 7163.2002
$ws.Range("M132").Value = -5903.8181
$ws.Range("N132").Value = -12223.2002
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2552.6155
$ws.Range("I7").Value = 2314.8572
$ws.Range("J7").Value = 2830
$ws.Range("K7").Value = 2314.8572
$ws.Range("L7").Value = 2830
$ws.Range("M7").Value = -2202.8572
$ws.Range("N7").Value = -3054
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""
$ws.Range("H46").Value = 696.6667
$ws.Range("I46").Value = 380
$ws.Range("J46").Value = 1140
$ws.Range("K46").Value = 380
$ws.Range("L46").Value = 1140
$ws.Range("M46").Value = -192
$ws.Range("N46").Value = -1516
$ws.Range("H126").Value = 2552.6155
$ws.Range("I126").Value = 2314.8572
$ws.Range("J126").Value = 2830
$ws.Range("K126").Value = 6944.571599999999
$ws.Range("L126").Value = 8490
$ws.Range("M126").Value = -4474.571599999999
$ws.Range("N126").Value = -13430
$ws.Range("H132").Value = 3476.2856
$ws.Range("I132").Value = 3050.2144
$ws.Range("J132").Value = 4328.4287
$ws.Range("K132").Value = 9150.643199999999
$ws.Range("L132").Value = 12985.2861
$ws.Range("M132").Value = -6620.643199999999
$ws.Range("N132").Value = -18045.2861
$ws.Range("H136").Value = 1944.75
$ws.Range("I136").Value = 1454.6923
$ws.Range("J136").Value = 4068.3333
$ws.Range("K136").Value = 4364.0769
$ws.Range("L136").Value = 12204.9999
$ws.Range("M136").Value = -1814.0769
$ws.Range("N136").Value = -17304.9999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 986.6875
$ws.Range("I107").Value = 848.9
$ws.Range("J107").Value = 1216.3334
$ws.Range("K107").Value = 2546.7
$ws.Range("L107").Value = 3649.0002
$ws.Range("M107").Value = -626.6999999999998
$ws.Range("N107").Value = -7489.0002
$ws.Range("H132").Value = 6077870
$ws.Range("I132").Value = 1345.65
$ws.Range("K132").Value = 4036.95
$ws.Range("M132").Value = -1506.95
$ws.Range("H136").Value = 2549.4
$ws.Range("I136").Value = 2439.1538
$ws.Range("K136").Value = 7317.4614
$ws.Range("M136").Value = -4767.4614
